# updated links and dates
# Slide 20, shape "CustomShape 2": change the registration-date aside on the
# "Advanced R course for data analysis and visualisation" line from
# "(registration open 2nd Feb)" to "(Starts 18th May)".

$p = $ppt.ActivePresentation

$targetSlideIndex = 0
$targetShapeIndex = 0
$targetParaIndex = 0

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if (-not $shape.HasTextFrame) { continue }
        $tf = $shape.TextFrame
        if (-not $tf.HasText) { continue }
        $tr = $tf.TextRange
        $paraCount = $tr.Paragraphs().Count
        for ($pi = 1; $pi -le $paraCount; $pi++) {
            $para = $tr.Paragraphs($pi, 1)
            if ($para.Text -like "*registration open*") {
                $targetSlideIndex = $si
                $targetShapeIndex = $shi
                $targetParaIndex = $pi
            }
        }
    }
}

if ($targetSlideIndex -eq 0) {
    throw "Could not locate paragraph containing 'registration open'"
}

$slide = $p.Slides.Item($targetSlideIndex)
$shape = $slide.Shapes.Item($targetShapeIndex)
$para = $shape.TextFrame.TextRange.Paragraphs($targetParaIndex, 1)

$runCount = $para.Runs().Count
for ($ri = 1; $ri -le $runCount; $ri++) {
    $run = $para.Runs($ri, 1)
    $txt = $run.Text
    if ($txt -eq "(registration open 2") {
        $run.Text = "(Starts 18"
    } elseif ($txt -eq "nd") {
        $run.Text = "th"
    } elseif ($txt -eq " Feb)") {
        $run.Text = " May)"
    }
}

Write-Output "Updated paragraph: [$($shape.TextFrame.TextRange.Paragraphs($targetParaIndex, 1).Text)]"
